# Add a new treatment (treatment_id 15, "calendarBar" / no interaction)
# to the "Web Parameters" sheet, mirroring the existing treatment_id 14
# ("calendarWord" / no interaction) rows, with the comment text updated to
# "Calendar year view with bar and no interaction."

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Web Parameters")
$ws.Activate()

$newComment = "Calendar year view with bar and no interaction."

$data = @(
    @{ row = 26; A = 15; B = 1; F = 300; H = 44593; I = 700;  K = 44703; L = 1100 },
    @{ row = 27; A = 15; B = 2; F = 500; H = 44621; I = 800;  K = 44724; L = 1100 },
    @{ row = 28; A = 15; B = 3; F = 300; H = 44652; I = 1000; K = 44757; L = 1100 }
)

foreach ($entry in $data) {
    $r = $entry.row

    $ws.Cells.Item($r, 1).Value = $entry.A      # A: treatment_id
    $ws.Cells.Item($r, 2).Value = $entry.B      # B: titration series no.
    $ws.Cells.Item($r, 3).Value = "calendarBar" # C: view_type
    $ws.Cells.Item($r, 4).Value = "none"        # D: interaction
    $ws.Cells.Item($r, 5).Value = "none"        # E: variable_amount
    $ws.Cells.Item($r, 6).Value = $entry.F      # F: amount_earlier

    $ws.Cells.Item($r, 8).Value = $entry.H      # H: date_earlier
    $ws.Cells.Item($r, 8).NumberFormat = "m/d/yyyy;@"

    $ws.Cells.Item($r, 9).Value = $entry.I      # I: amount_later

    $ws.Cells.Item($r, 11).Value = $entry.K     # K: date_later
    $ws.Cells.Item($r, 11).NumberFormat = "m/d/yyyy;@"

    $ws.Cells.Item($r, 12).Value = $entry.L     # L: max_amount

    $ws.Cells.Item($r, 14).Value = 100          # N: horizontal_pixels
    $ws.Cells.Item($r, 15).Value = 100          # O: vertical_pixels

    $ws.Cells.Item($r, 20).Value = 10           # T: width_in
    $ws.Cells.Item($r, 21).Value = 8            # U: height_in

    $ws.Cells.Item($r, 22).Value = $newComment  # V: comment
}

$ws.Range("E26").Select()
